$wb = $excel.ActiveWorkbook

# --- Sheet "ENGLISH": append 5 new vocabulary rows (182-186) ---
$ws1 = $wb.Worksheets.Item("ENGLISH")

$english = @(
    @("astounding", "", "amazing;impressive;notable", "", 0, "2021-11-23 13:31:52.160918", "", ""),
    @("tangible", "", "real;touchable", "", 0, "2021-11-23 13:32:51.714678", "", ""),
    @("growl", "", "snarl;say roughly", "", 0, "2021-11-23 13:38:21.392315", "", ""),
    @("impose", "force (an unwelcome decision or ruling) on someone", "", "", 0, "2021-11-23 13:40:32.909102", "", ""),
    @("revolt", "", "rebel", "", 0, "2021-11-23 13:41:24.307841", "", "")
)

$startRow = 182
for ($i = 0; $i -lt $english.Length; $i++) {
    $row = $startRow + $i
    $data = $english[$i]
    for ($col = 1; $col -le 8; $col++) {
        $ws1.Cells.Item($row, $col).Value = $data[$col - 1]
    }
}

# --- Sheet "NOTES": append 7 new note rows (49-55) ---
$ws2 = $wb.Worksheets.Item("NOTES")

$notes = @(
    @("Same habits, same results", "habit"),
    @("Changes that seem small and unimportant at first will compound into remarkable results if you are willing to stick with them for years", "habit"),
    @("In the long run, the quality of our lives often depends on the quality of our habits", "habit"),
    @("Mastery requires patience", ""),
    @("Goals are good for setting a direction, but systems are best for making progress", ""),
    @("True behavior change is identity change", "habit"),
    @("Habits create freedom", "habit")
)

$startRow2 = 49
for ($i = 0; $i -lt $notes.Length; $i++) {
    $row = $startRow2 + $i
    $data = $notes[$i]
    $ws2.Cells.Item($row, 1).Value = $data[0]
    $ws2.Cells.Item($row, 2).Value = $data[1]
}
